# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Mapping of sheet -> cell -> new value, derived from the OOXML diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1805
$ws1.Range("F9").Value = 1768
$ws1.Range("F15").Value = 12932
$ws1.Range("F18").Value = 750

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 179

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 179
$ws4.Range("F5").Value = 1805
$ws4.Range("F14").Value = 1768
$ws4.Range("F21").Value = 12932
$ws4.Range("F24").Value = 750
